# Auto-generated edit script: apply cached-value updates to Jenova_Profits sheets
# (values refreshed by the scheduled market-data runner; cells hold static
#  numbers, not formulas, so each changed cell is written directly.)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 191.5
$ws.Range("I33").Value = 191.5
$ws.Range("K33").Value = 191.5
$ws.Range("M33").Value = 37.5
$ws.Range("H64").Value = 6588.2354
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 6588.2354
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H98").Value = 1148.9565
$ws.Range("I98").Value = 696.35
$ws.Range("K98").Value = 696.35
$ws.Range("M98").Value = 801.65
$ws.Range("H113").Value = 2154.75
$ws.Range("I113").Value = 2039.6666
$ws.Range("K113").Value = 2039.6666
$ws.Range("M113").Value = 1214.3334
$ws.Range("H116").Value = 11955.866
$ws.Range("I116").Value = 5665.4287
$ws.Range("K116").Value = 5665.4287
$ws.Range("M116").Value = -2223.4287
$ws.Range("H122").Value = 1148.9565
$ws.Range("I122").Value = 696.35
$ws.Range("K122").Value = 2089.05
$ws.Range("M122").Value = 360.9499999999998
$ws.Range("H125").Value = 18520636
$ws.Range("I125").Value = 1067
$ws.Range("K125").Value = 9603
$ws.Range("M125").Value = -7143
$ws.Range("H129").Value = 1596.6552
$ws.Range("I129").Value = 633.7273
$ws.Range("J129").Value = 2185.111
$ws.Range("K129").Value = 1901.1819
$ws.Range("L129").Value = 6555.333
$ws.Range("M129").Value = 3098.8181
$ws.Range("N129").Value = -16555.333
$ws.Range("H132").Value = 4630
$ws.Range("I132").Value = 4739.4546
$ws.Range("K132").Value = 14218.3638
$ws.Range("M132").Value = -11688.3638
$ws.Range("H137").Value = 3887.6667
$ws.Range("I137").Value = 3212.7144
$ws.Range("K137").Value = 9638.143199999999
$ws.Range("M137").Value = -7088.143199999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4672.8237
$ws.Range("I61").Value = 4429.2666
$ws.Range("K61").Value = 4429.2666
$ws.Range("M61").Value = -4217.2666
$ws.Range("H74").Value = 2346.4614
$ws.Range("I74").Value = 2813
$ws.Range("K74").Value = 2813
$ws.Range("M74").Value = -1939
$ws.Range("H77").Value = 2346.4614
$ws.Range("I77").Value = 2813
$ws.Range("K77").Value = 14065
$ws.Range("M77").Value = -9697
$ws.Range("H102").Value = 1783.5
$ws.Range("I102").Value = 1354.7273
$ws.Range("K102").Value = 1354.7273
$ws.Range("M102").Value = 267.2727
$ws.Range("H110").Value = 502893.1
$ws.Range("I110").Value = 628291.4
$ws.Range("J110").Value = 1300
$ws.Range("K110").Value = 628291.4
$ws.Range("L110").Value = 1300
$ws.Range("M110").Value = -626246.4
$ws.Range("N110").Value = -5390
$ws.Range("H136").Value = 4672.8237
$ws.Range("I136").Value = 4429.2666
$ws.Range("K136").Value = 13287.7998
$ws.Range("M136").Value = -10737.7998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 631643.5
$ws.Range("I86").Value = 851783
$ws.Range("K86").Value = 851783
$ws.Range("M86").Value = -850660
$ws.Range("H89").Value = 631643.5
$ws.Range("I89").Value = 851783
$ws.Range("K89").Value = 4258915
$ws.Range("M89").Value = -4253299
$ws.Range("H105").Value = 3415.8572
$ws.Range("I105").Value = 3415.8572
$ws.Range("K105").Value = 3415.8572
$ws.Range("M105").Value = -1668.8572

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 10165
$ws.Range("I16").Value = 3663
$ws.Range("J16").Value = 29671
$ws.Range("K16").Value = 3663
$ws.Range("L16").Value = 29671
$ws.Range("M16").Value = -3376
$ws.Range("N16").Value = -30245
$ws.Range("H58").Value = 8015.7144
$ws.Range("I58").Value = 5222
$ws.Range("J58").Value = 15000
$ws.Range("K58").Value = 5222
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -5019
$ws.Range("N58").Value = -15406
$ws.Range("H94").Value = 1016.1539
$ws.Range("I94").Value = 503
$ws.Range("J94").Value = 1244.2222
$ws.Range("K94").Value = 503
$ws.Range("L94").Value = 1244.2222
$ws.Range("M94").Value = -52
$ws.Range("N94").Value = -2146.2222
$ws.Range("H98").Value = 60840.5
$ws.Range("J98").Value = 60840.5
$ws.Range("L98").Value = 60840.5
$ws.Range("N98").Value = -65332.5
$ws.Range("H99").Value = 4313
$ws.Range("I99").Value = 4126.5
$ws.Range("K99").Value = 4126.5
$ws.Range("M99").Value = -2628.5
$ws.Range("H113").Value = 10165
$ws.Range("I113").Value = 3663
$ws.Range("J113").Value = 29671
$ws.Range("K113").Value = 3663
$ws.Range("L113").Value = 29671
$ws.Range("M113").Value = -1493
$ws.Range("N113").Value = -34011
$ws.Range("H122").Value = 3238.611
$ws.Range("I122").Value = 2559.3845
$ws.Range("J122").Value = 5004.6
$ws.Range("K122").Value = 7678.1535
$ws.Range("L122").Value = 15013.8
$ws.Range("M122").Value = -5228.1535
$ws.Range("N122").Value = -19913.8
$ws.Range("H126").Value = 4313
$ws.Range("I126").Value = 4126.5
$ws.Range("K126").Value = 12379.5
$ws.Range("M126").Value = -9909.5
$ws.Range("H132").Value = 2747.8823
$ws.Range("J132").Value = 3145.1428
$ws.Range("L132").Value = 9435.428400000001
$ws.Range("N132").Value = -14495.4284
$ws.Range("H134").Value = 529964.0600000001
$ws.Range("I134").Value = 3779.2942
$ws.Range("K134").Value = 11337.8826
$ws.Range("M134").Value = -8802.882599999999
$ws.Range("H136").Value = 8015.7144
$ws.Range("I136").Value = 5222
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 15666
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -13116
$ws.Range("N136").Value = -50100
$ws.Range("H141").Value = 405138.7
$ws.Range("J141").Value = 479474.38
$ws.Range("L141").Value = 479474.38
$ws.Range("N141").Value = -489834.38

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 27
$ws.Range("I19").Value = 5
$ws.Range("K19").Value = 15
$ws.Range("M19").Value = 159

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1642.7715
$ws.Range("I102").Value = 1089.5862
$ws.Range("K102").Value = 1089.5862
$ws.Range("M102").Value = 532.4138
$ws.Range("H113").Value = 669788.3
$ws.Range("I113").Value = 1430189.4
$ws.Range("K113").Value = 1430189.4
$ws.Range("M113").Value = -1428019.4
$ws.Range("H122").Value = 3877.7144
$ws.Range("I122").Value = 1298.5
$ws.Range("J122").Value = 7316.6665
$ws.Range("K122").Value = 3895.5
$ws.Range("L122").Value = 21949.9995
$ws.Range("M122").Value = -1445.5
$ws.Range("N122").Value = -26849.9995
$ws.Range("H126").Value = 4571.4287
$ws.Range("H132").Value = 112234.4
$ws.Range("I132").Value = 13593.777
$ws.Range("K132").Value = 40781.331
$ws.Range("M132").Value = -38251.331

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3023.3225
$ws.Range("I40").Value = 3074.1
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 3074.1
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -2938.1
$ws.Range("N40").Value = -1772
$ws.Range("I122").Value = 5999
$ws.Range("K122").Value = 17997
$ws.Range("M122").Value = -15547
$ws.Range("H132").Value = 3974.375
$ws.Range("I132").Value = 1299.3334
$ws.Range("J132").Value = 11999.5
$ws.Range("K132").Value = 3898.0002
$ws.Range("L132").Value = 35998.5
$ws.Range("M132").Value = -1368.0002
$ws.Range("N132").Value = -41058.5
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 90000
$ws.Range("I64").Value = 90000
$ws.Range("K64").Value = 90000
$ws.Range("M64").Value = -89752
$ws.Range("H67").Value = 90000
$ws.Range("I67").Value = 90000
$ws.Range("K67").Value = 90000
$ws.Range("M67").Value = -89142
$ws.Range("H81").Value = 4617.9443
$ws.Range("I81").Value = 2194.9167
$ws.Range("J81").Value = 9464
$ws.Range("K81").Value = 4389.8334
$ws.Range("L81").Value = 18928
$ws.Range("M81").Value = -3328.8334
$ws.Range("N81").Value = -21050
$ws.Range("H84").Value = 4617.9443
$ws.Range("I84").Value = 2194.9167
$ws.Range("J84").Value = 9464
$ws.Range("K84").Value = 21949.167
$ws.Range("L84").Value = 94640
$ws.Range("M84").Value = -16645.167
$ws.Range("N84").Value = -105248
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 38154.83
$ws.Range("I132").Value = 2520.3157
$ws.Range("J132").Value = 105860.4
$ws.Range("K132").Value = 7560.9471
$ws.Range("L132").Value = 317581.2
$ws.Range("M132").Value = -5030.9471
$ws.Range("N132").Value = -322641.2
$ws.Range("H136").Value = 117671.57
$ws.Range("I136").Value = 25811.389
$ws.Range("J136").Value = 668832.7
$ws.Range("K136").Value = 77434.167
$ws.Range("L136").Value = 2006498.1
$ws.Range("M136").Value = -74884.167
$ws.Range("N136").Value = -2011598.1
